# Insert a new data row for "Macroferia Regional de Talca - Papa" at row 591,
# pushing all existing rows from 591..657 down to 592..658.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 591 (shifts 591:657 -> 592:658)
$ws.Rows.Item(591).Insert()

# Populate the new row 591 with the new record's data
$ws.Cells.Item(591, 1).Value = 5
$ws.Cells.Item(591, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(591, 3).Value = "Maule"
$ws.Cells.Item(591, 4).Value2 = 44918
$ws.Cells.Item(591, 5).Value = 7
$ws.Cells.Item(591, 6).Value = 100114001
$ws.Cells.Item(591, 7).Value = "Papa"
$ws.Cells.Item(591, 8).Value = "Rodeo"
$ws.Cells.Item(591, 9).Value = "1a (cosecha)"
$ws.Cells.Item(591, 10).Value = 1500
$ws.Cells.Item(591, 11).Value = 12000
$ws.Cells.Item(591, 12).Value = 12000
$ws.Cells.Item(591, 13).Value = 12000
$ws.Cells.Item(591, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(591, 15).Value = "Región del Maule"
$ws.Cells.Item(591, 16).Value = 480
$ws.Cells.Item(591, 17).Value = 25
$ws.Cells.Item(591, 18).Value = "Hortaliza"
